{"js": "// Office.js (Word JavaScript API) edit script.\n// The document's title date and every \"##\u00f7#=\" division-problem cell get\n// replaced with new values. All text runs change, one-for-one, in\n// document order (title paragraph first, then each table row's cells\n// left-to-right, top-to-bottom) - no rows/cells are added or removed,\n// so we can do a single positional pass over context.document.body.paragraphs.\n\nconst oldTexts = [\n  \"2026-01-12 Monday\",\n  \"68\u00f77=\", \"19\u00f78=\", \"67\u00f75=\", \"71\u00f77=\", \"14\u00f77=\",\n  \"64\u00f73=\", \"30\u00f79=\", \"26\u00f76=\", \"94\u00f78=\", \"82\u00f75=\",\n  \"36\u00f77=\", \"53\u00f72=\", \"47\u00f72=\", \"39\u00f79=\", \"10\u00f73=\",\n  \"31\u00f72=\", \"93\u00f76=\", \"28\u00f79=\", \"65\u00f79=\", \"13\u00f74=\",\n  \"50\u00f73=\", \"39\u00f79=\", \"54\u00f74=\", \"17\u00f73=\", \"34\u00f73=\"\n];\n\nconst newTexts = [\n  \"2026-01-13 Tuesday\",\n  \"51\u00f78=\", \"96\u00f78=\", \"63\u00f72=\", \"96\u00f78=\", \"45\u00f78=\",\n  \"72\u00f79=\", \"61\u00f78=\", \"57\u00f72=\", \"48\u00f79=\", \"73\u00f75=\",\n  \"92\u00f77=\", \"46\u00f78=\", \"84\u00f74=\", \"81\u00f74=\", \"37\u00f75=\",\n  \"46\u00f79=\", \"75\u00f74=\", \"50\u00f78=\", \"51\u00f73=\", \"68\u00f76=\",\n  \"16\u00f79=\", \"35\u00f73=\", \"39\u00f79=\", \"26\u00f76=\", \"89\u00f78=\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < oldTexts.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === oldTexts[idx]) {\n    para.insertText(newTexts[idx], \"Replace\");\n    idx++;\n  }\n}\nawait context.sync();\n\nif (idx !== oldTexts.length) {\n  throw new Error(\n    \"Expected to replace \" + oldTexts.length + \" paragraphs but only matched \" + idx\n  );\n}\n", "ps1": "# Word COM interop edit script.\n# The document's title date and every \"##\u00f7#=\" division-problem cell get\n# replaced with new values. All text runs change, one-for-one, in\n# document order (title paragraph first, then each table row's cells\n# left-to-right, top-to-bottom) - no rows/cells are added or removed,\n# so we can do a single positional pass over $d.Paragraphs.\n\n$d = $word.ActiveDocument\n\n$oldTexts = @(\n  \"2026-01-12 Monday\",\n  \"68\u00f77=\", \"19\u00f78=\", \"67\u00f75=\", \"71\u00f77=\", \"14\u00f77=\",\n  \"64\u00f73=\", \"30\u00f79=\", \"26\u00f76=\", \"94\u00f78=\", \"82\u00f75=\",\n  \"36\u00f77=\", \"53\u00f72=\", \"47\u00f72=\", \"39\u00f79=\", \"10\u00f73=\",\n  \"31\u00f72=\", \"93\u00f76=\", \"28\u00f79=\", \"65\u00f79=\", \"13\u00f74=\",\n  \"50\u00f73=\", \"39\u00f79=\", \"54\u00f74=\", \"17\u00f73=\", \"34\u00f73=\"\n)\n\n$newTexts = @(\n  \"2026-01-13 Tuesday\",\n  \"51\u00f78=\", \"96\u00f78=\", \"63\u00f72=\", \"96\u00f78=\", \"45\u00f78=\",\n  \"72\u00f79=\", \"61\u00f78=\", \"57\u00f72=\", \"48\u00f79=\", \"73\u00f75=\",\n  \"92\u00f77=\", \"46\u00f78=\", \"84\u00f74=\", \"81\u00f74=\", \"37\u00f75=\",\n  \"46\u00f79=\", \"75\u00f74=\", \"50\u00f78=\", \"51\u00f73=\", \"68\u00f76=\",\n  \"16\u00f79=\", \"35\u00f73=\", \"39\u00f79=\", \"26\u00f76=\", \"89\u00f78=\"\n)\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    if ($idx -ge $oldTexts.Length) { break }\n    # Range.Text includes the trailing paragraph mark (CR); strip it before comparing.\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $oldTexts[$idx]) {\n        $p.Range.Text = $newTexts[$idx]\n        $idx = $idx + 1\n    }\n}\n\nif ($idx -ne $oldTexts.Length) {\n    throw \"Expected to replace $($oldTexts.Length) paragraphs but only matched $idx\"\n}\n\nWrite-Output \"replaced: $idx\"\n"}
